# Katalon API Project commit
# Update the test credentials stored in row 2 (A2/B2) of the data file
# with the latest generated username/password pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test18042025@821.com"
$ws.Range("B2").Value = "Test18042025@821"
